$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update summary figures
$ws.Range("E11").Value = 108560      # Valor Mora total
$ws.Range("C13").Value = 1           # Cant. Trabajadores
$ws.Range("F13").Value = 2           # Cant. Periodos

# Remove the old middle worker row (FABIOLA MORENO ZABALETA / 45360660).
# Deleting row 17 shifts row 18 (PPT / JACKELINE VILLANUEVA) up into row 17,
# preserving that row's own (bottom-border) styling.
$ws.Rows(17).Delete()

# Overwrite the row that shifted up (now row 17) with the second period's
# data for the remaining worker (TEOFILO PEREZ FORTICHE), period 2508.
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73506856"
$ws.Range("D17").Value = "TEOFILO PEREZ FORTICHE"
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 54280
$ws.Range("G17").Value = 1357000

# Nudge the logo image left (it was manually repositioned), keeping its
# original size.
$logo = $ws.Shapes.Item(1)
$logo.Left = 53.59055118110236
$logo.Width = 76.81889763779527
$logo.Height = 48.188976377952756
